$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 1.37
$ws.Range("E3").Value = 1.29
$ws.Range("G3").Value = 0.61
$ws.Range("F4").Value = 1.1
$ws.Range("C5").Value = 1.35
$ws.Range("D6").Value = 1.55
$ws.Range("G6").Value = 1.05
$ws.Range("C7").Value = 2.22
$ws.Range("F7").Value = 1.44
